$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Insert a new worksheet named "Merge" right before "ToUnhide"
#    (so it lands as the 10th sheet, pushing ToUnhide to 11th).
# ------------------------------------------------------------------
$toUnhide = $wb.Worksheets.Item("ToUnhide")
$mergeSheet = $wb.Worksheets.Add($toUnhide)
$mergeSheet.Name = "Merge"

# Put the explanatory note in B2 and select C2, matching the authored sheet.
$mergeSheet.Range("B2").Value = "D4 to H8 should be merged…"
$mergeSheet.Columns("B").ColumnWidth = 26.7109375
$mergeSheet.Range("C2").Select()

# ------------------------------------------------------------------
# 2) Demo1: rows 4-9 shrink from height 24 to 23.25 (cosmetic re-fit
#    tied to the "should be merged" note being added to the workbook).
# ------------------------------------------------------------------
$demo1 = $wb.Worksheets.Item("Demo1")
$demo1.Rows("4:9").RowHeight = 23.25
